$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New B,C,D,E,G values per row (F column is unchanged; G = B+C+D+E)
$data = @{
    2  = @(0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1.35982162114495)
    3  = @(0.3048080303191223, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.274871460341982)
    4  = @(0.127881588408715, 0.00007097389502863649, 0.8054896365839992, 8.660232485948974, 9.593674684836717)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    6  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    7  = @(1.459612070389937, 3099.503889238888, 3.900430680208489, 645.3272768299601, 3750.191208819447)
    8  = @(1.459612070389937, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 6.169612356620297)
    9  = @(0.0008583669626518464, 0.00007097389502863649, 0.8054896365839992, 0.496779210170732, 1.303198187612412)
    10 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    11 = @(1.459612070389937, 3099.503889238888, 337.1190423067083, 8.660232485948974, 3446.742776101936)
    12 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 11.945164432584)
    13 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 31.61296591696135)
    14 = @(0.6753301551942219, 1.667794583268128, 337.1190423067083, 645.3272768299601, 984.7894438751307)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E
    $ws.Cells.Item($r, 7).Value = $vals[4]   # G
}
